$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 40 and 41: Hedera and MXToken swap places (name + link)
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"

# Update Price (D) and Volume(1h) (E) columns for rows 2-51.
# Price values that look like plain numbers must be forced back to
# text (quote-prefix) so Excel does not silently convert them to
# numbers, then the style is reset to Normal so no stray number
# format is left behind on the cell.
$ws.Range("D2").Value = "26.559.73"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.815.09"
$ws.Range("E3").Value = "  +0.65%  "
$c = $ws.Range("D4")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.79%  "
$c = $ws.Range("D5")
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "
$c = $ws.Range("D6")
$c.Value = "'305.88"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "
$c = $ws.Range("D7")
$c.Value = "'0.4530"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.60%  "
$c = $ws.Range("D8")
$c.Value = "'0.3611"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "
$c = $ws.Range("D9")
$c.Value = "'46.35"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.98%  "
$c = $ws.Range("D10")
$c.Value = "'0.07105"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.14%  "
$c = $ws.Range("D11")
$c.Value = "'0.8967"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.34%  "
$c = $ws.Range("D12")
$c.Value = "'0.07790"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.86%  "
$c = $ws.Range("D13")
$c.Value = "'19.42"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "1.789.75"
$ws.Range("E14").Value = "  -0.62%  "
$c = $ws.Range("D15")
$c.Value = "'5.288"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.51%  "
$c = $ws.Range("D16")
$c.Value = "'6.322"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "
$c = $ws.Range("D17")
$c.Value = "'85.36"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.60%  "
$c = $ws.Range("D18")
$c.Value = "'1.003"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.68%  "
$c = $ws.Range("D19")
$c.Value = "'0.000008614"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.57%  "
$c = $ws.Range("D20")
$c.Value = "'1.002"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "26.586.49"
$ws.Range("E21").Value = "  +0.32%  "
$c = $ws.Range("D22")
$c.Value = "'14.24"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$c = $ws.Range("D23")
$c.Value = "'4.974"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "
$c = $ws.Range("D24")
$c.Value = "'10.55"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").Value = "2.003.66"
$ws.Range("E25").Value = "  -2.27%  "
$c = $ws.Range("D26")
$c.Value = "'1.957"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.62%  "
$c = $ws.Range("D27")
$c.Value = "'150.86"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "
$c = $ws.Range("D28")
$c.Value = "'17.84"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.36%  "
$c = $ws.Range("D29")
$c.Value = "'2.060"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.74%  "
$c = $ws.Range("D30")
$c.Value = "'112.41"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.04%  "
$c = $ws.Range("D31")
$c.Value = "'4.860"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.99%  "
$c = $ws.Range("D32")
$c.Value = "'0.08703"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "
$c = $ws.Range("D33")
$c.Value = "'3.126"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.64%  "
$c = $ws.Range("D34")
$c.Value = "'0.7494"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.06%  "
$c = $ws.Range("D35")
$c.Value = "'2.768"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +9.33%  "
$c = $ws.Range("D36")
$c.Value = "'4.451"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.54%  "
$c = $ws.Range("D37")
$c.Value = "'1.114"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "
$c = $ws.Range("D38")
$c.Value = "'1.071"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.60%  "
$c = $ws.Range("D39")
$c.Value = "'0.01932"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.22%  "
$c = $ws.Range("D40")
$c.Value = "'0.05119"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "
$c = $ws.Range("D41")
$c.Value = "'2.899"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "
$c = $ws.Range("D42")
$c.Value = "'0.5097"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.55%  "
$c = $ws.Range("D43")
$c.Value = "'6.744"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "
$c = $ws.Range("D44")
$c.Value = "'0.1510"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.32%  "
$c = $ws.Range("D45")
$c.Value = "'8.075"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "
$c = $ws.Range("D46")
$c.Value = "'0.4745"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.51%  "
$c = $ws.Range("D47")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "
$c = $ws.Range("D48")
$c.Value = "'10.03"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.13%  "
$c = $ws.Range("D49")
$c.Value = "'100.85"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.56%  "
$c = $ws.Range("D50")
$c.Value = "'1.582"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "
$c = $ws.Range("D51")
$c.Value = "'0.05982"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
